$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.001" or
# "0.06700" keep their exact literal formatting instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.439.19"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "1.893.31"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "238.18"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.4901"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "0.2939"
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").Value = "0.06700"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Value = "1.877.41"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").Value = "17.04"
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").Value = "0.07347"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "5.142"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").Value = "88.06"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "0.6654"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "30.419.01"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "13.47"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "0.000007825"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "2.143.06"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "5.328"
$ws.Range("E21").Value = "  +12.14%  "
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "190.13"
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("D24").Value = "6.144"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "9.507"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "164.32"
$ws.Range("E26").Value = "  +3.03%  "
$ws.Range("D27").Value = "18.27"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  +5.72%  "
$ws.Range("D29").Value = "1.465"
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("D30").Value = "4.356"
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "0.09163"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("D32").Value = "4.066"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").Value = "0.05215"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "0.7420"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").Value = "1.101"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").Value = "2.716"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "0.01815"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").Value = "2.673"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "0.9220"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "2.037"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").Value = "0.4411"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "5.946"
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("D43").Value = "106.35"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").Value = "0.9937"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "68.68"
$ws.Range("E45").Value = "  +19.84%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.1378"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").Value = "7.590"
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("D48").Value = "8.977"
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("D49").Value = "34.98"
$ws.Range("E49").Value = "  +5.07%  "
$ws.Range("D50").Value = "0.05821"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").Value = "0.3952"
$ws.Range("E51").Value = "  -5.14%  "
